$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Align the duplicated "Center A Ben Mansour" rows (3 and 4) with row 2's
# latitude/longitude so they share the same coordinates and shared-string text.
# Using Copy (instead of .Value) preserves text typing / default styling exactly
# like the source cell (avoids Excel auto-coercing the leading-space text to a number).
$ws.Range("G2").Copy($ws.Range("G3"))
$ws.Range("G2").Copy($ws.Range("G4"))
$ws.Range("H2").Copy($ws.Range("H3"))
$ws.Range("H2").Copy($ws.Range("H4"))

# Two new Mac-Addresses (kiosks) added to every registration center:
# rows 2-4 (Center A Ben Mansour, 3 duplicate rows) now have 3 kiosks,
# all remaining centers (rows 5-46) now have 2 kiosks instead of 1.
$ws.Range("L2").Value = 3
$ws.Range("L3").Value = 3
$ws.Range("L4").Value = 3

for ($row = 5; $row -le 46; $row++) {
    $ws.Cells.Item($row, 12).Value = 2
}
